$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '91.808.70'
$ws.Range("E2").Value = '  -2.60%  '

$ws.Range("D3").Value = '3.327.33'
$ws.Range("E3").Value = '  -4.03%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.63'
$ws.Range("E5").Value = '  -2.94%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '614.70'
$ws.Range("E6").Value = '  -4.05%  '

$ws.Range("E7").Value = '  -1.90%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.383'
$ws.Range("E8").Value = '  -3.25%  '

$ws.Range("E9").Value = '  +0.06%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.955'
$ws.Range("E10").Value = '  -0.47%  '

$ws.Range("D11").Value = '3.327.18'
$ws.Range("E11").Value = '  -3.98%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.93'
$ws.Range("E12").Value = '  +1.55%  '

$ws.Range("E13").Value = '  -1.69%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.14'
$ws.Range("E14").Value = '  +0.28%  '

$ws.Range("D15").Value = '91.700.09'
$ws.Range("E15").Value = '  -2.51%  '

$ws.Range("D16").Value = '3.944.58'
$ws.Range("E16").Value = '  -4.12%  '

$ws.Range("E17").Value = '  -3.26%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.07'
$ws.Range("E18").Value = '  -3.33%  '

$ws.Range("D19").Value = '3.329.58'
$ws.Range("E19").Value = '  -4.05%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.32'
$ws.Range("E20").Value = '  -2.27%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.82'
$ws.Range("E21").Value = '  -4.50%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.41'
$ws.Range("E22").Value = '  +5.33%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '492.61'
$ws.Range("E23").Value = '  -1.07%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.442'
$ws.Range("E24").Value = '  -9.84%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.43'
$ws.Range("E25").Value = '  -1.63%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000182'
$ws.Range("E26").Value = '  -4.48%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '89.92'
$ws.Range("E27").Value = '  -0.95%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.90'
$ws.Range("E28").Value = '  -0.74%  '

$ws.Range("D29").Value = '3.500.60'
$ws.Range("E29").Value = '  -3.87%  '

$ws.Range("E30").Value = '  +0.02%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.12'
$ws.Range("E31").Value = '  -5.18%  '

$ws.Range("E32").Value = '  +1.95%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.59'
$ws.Range("E33").Value = '  -4.87%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.02'
$ws.Range("E34").Value = '  +2.10%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.172'
$ws.Range("E35").Value = '  -4.60%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '28.32'
$ws.Range("E36").Value = '  -6.48%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.526'
$ws.Range("E37").Value = '  -5.28%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '565.82'
$ws.Range("E38").Value = '  +1.85%  '

$ws.Range("B39").Value = 'USDe'
$ws.Range("C39").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  -0.05%  '

$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.36'
$ws.Range("E40").Value = '  -3.65%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.148'
$ws.Range("E41").Value = '  -1.43%  '

$ws.Range("E42").Value = '  -5.59%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.867'
$ws.Range("E43").Value = '  -6.99%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '23.66'
$ws.Range("E44").Value = '  -1.62%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.67'
$ws.Range("E45").Value = '  -2.49%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0412'
$ws.Range("E46").Value = '  +0.27%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.60'
$ws.Range("E47").Value = '  +2.16%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.40'
$ws.Range("E48").Value = '  -2.41%  '

$ws.Range("E49").Value = '  -1.78%  '

$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '51.77'
$ws.Range("E50").Value = '  -2.33%  '

$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.96'
$ws.Range("E51").Value = '  -0.67%  '
